$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Fill in the three previously-empty assignment blocks on row 2 with the
#    new "4 UE" / "5 UE" / "6 UE" grades (all scored 100, "Good Job!").
#    Order matters so the shared-string table ends up built in the same
#    sequence as the authored workbook.
# ---------------------------------------------------------------------------
$ws.Range("R2").Value = "4 UE"
$ws.Range("S2").Value = 100
$ws.Range("T2").Value = "Good Job!"

$ws.Range("V2").Value = "5 UE"
$ws.Range("W2").Value = 100
$ws.Range("X2").Value = "Good Job!"

$ws.Range("Z2").Value = "6 UE"
$ws.Range("AA2").Value = 100
$ws.Range("AB2").Value = "Good Job!"

# ---------------------------------------------------------------------------
# 2) Update existing grade data for student "2 UE" (row 2): grade 10 -> 100,
#    and the comment changes from the old sphere-comment to "I Saw the fix".
# ---------------------------------------------------------------------------
$ws.Range("K2").Value = 100
$ws.Range("L2").Value = "I Saw the fix"

# ---------------------------------------------------------------------------
# 3) Add three more assignment-block columns (AC:AN) mirroring the existing
#    "Assignment / Grade / Comments / (spacer)" pattern, headers on row 1.
# ---------------------------------------------------------------------------
$ws.Range("AD1").Value = "Assignment"
$ws.Range("AE1").Value = "Grade"
$ws.Range("AF1").Value = "Comments"

$ws.Range("AH1").Value = "Assignment"
$ws.Range("AI1").Value = "Grade"
$ws.Range("AJ1").Value = "Comments"

$ws.Range("AL1").Value = "Assignment"
$ws.Range("AM1").Value = "Grade"
$ws.Range("AN1").Value = "Comments"

# Copy formatting (styles only) from the last existing block (Y1:AB18, which
# covers header row 1 + all data rows 2-18) into each of the three new
# blocks so borders/fills/alignment match the rest of the table.
$ws.Range("Y1:AB18").Copy() | Out-Null
$ws.Range("AC1").PasteSpecial(-4122) | Out-Null
$ws.Range("Y1:AB18").Copy() | Out-Null
$ws.Range("AG1").PasteSpecial(-4122) | Out-Null
$ws.Range("Y1:AB18").Copy() | Out-Null
$ws.Range("AK1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Column widths for the new block (AC:AN)
$ws.Range("AC1:AN1").ColumnWidth = 15

# Row 1 is taller now that headers wrap across more columns
$ws.Rows.Item(1).RowHeight = 30

# ---------------------------------------------------------------------------
# 4) View state: the author scrolled right and selected L2.
# ---------------------------------------------------------------------------
$ws.Range("L2").Select() | Out-Null
